$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 23 ("https://www.instiglio.org/en/ve-dib/") which is a duplicate
# of row 11. This shifts rows 24 and 25 up by one (becoming 23 and 24).
$ws.Rows("23").Delete()
